$wb = $excel.ActiveWorkbook

# --- "Overview" sheet: row 3 is the b.md entry ----------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-20 00:43:25"

# --- "zh-cn" sheet: row 3 is the b.md source file --------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"

# "False" (and "True") look like booleans to the Value setter, so write it
# as a formula first, then paste-special as values into the real cell; that
# keeps the cell a plain shared-string "False" instead of a Boolean FALSE.
$zhcn.Range("ZZ1").Formula = '="False"'
$zhcn.Range("ZZ1").Copy()
$zhcn.Range("F3").PasteSpecial(-4163)
$zhcn.Range("ZZ1").Clear()

$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-20 00:43:20"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1e5c5a56ede401f5843448297f3362dfe22155a5/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/02f3e0b9a0482c7b2bddc6835acfa0e9be1c1c95/e2e/b.md."
$zhcn.Columns.Item(16).ColumnWidth = 40 - 5/6

# --- "de-de" sheet: row 3 is the b.md source file ---------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"

$dede.Range("ZZ1").Formula = '="False"'
$dede.Range("ZZ1").Copy()
$dede.Range("F3").PasteSpecial(-4163)
$dede.Range("ZZ1").Clear()

$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-08-20 00:43:25"
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1e5c5a56ede401f5843448297f3362dfe22155a5/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/02f3e0b9a0482c7b2bddc6835acfa0e9be1c1c95/e2e/b.md."
$dede.Columns.Item(16).ColumnWidth = 40 - 5/6
